$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Loc"
$ws.Range("B2").Value = "fff"
$ws.Range("C2").Value = "Loc.png"
$ws.Range("D2").Value = "ff"
$ws.Range("E2").Value = 0

$ws.Range("A2:E2").Select()
